$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "ECs", "Thbs2", "Itga6", "ECs", 3, 1, 1.924015333333333, 5.772046, 0.009877822204539638, 0.009877822204539637, 3, 1, 192.8285726666667, 578.485718, 0.7801188850698786, 0.7801188850698786, 371.0051305154475, 3339.046174639028, 0.007705875645123953, 0.007705875645123952),
    @(3, "ECs", "Thbs2", "Itga6", "FAPs", 3, 1, 1.924015333333333, 5.772046, 0.009877822204539638, 0.009877822204539637, 3, 1, 0.4209206666666667, 1.262762, 0.001702901995496819, 0.001702901995496819, 0.8098578167835556, 7.288720351052, 0.00001682096314327333858810315, 0.000016820963143273328423707783),
    @(4, "ECs", "Thbs2", "Itga6", "M2", 3, 1, 1.924015333333333, 5.772046, 0.009877822204539638, 0.009877822204539637, 3, 1, 45.70525533333333, 137.115766, 0.184907933193646, 0.184907933193646, 87.93761207524844, 791.438508677236, 0.001826487688295729, 0.001826487688295728),
    @(5, "ECs", "Thbs2", "Itga6", "sCs", 3, 1, 1.924015333333333, 5.772046, 0.009877822204539638, 0.009877822204539637, 3, 1, 8.223696, 24.671088, 0.0332702797409786, 0.0332702797409786, 15.822517200672, 142.402654806048, 0.0003286379079766838, 0.0003286379079766836),
    @(6, "FAPs", "Thbs2", "Itga6", "ECs", 3, 1, 188.2309416666667, 564.692825, 0.9663705600283187, 0.9663705600283184, 3, 1, 192.8285726666667, 578.485718, 0.7801188850698786, 0.7801188850698786, 36296.30381328593, 326666.7343195733, 0.7538839238536461, 0.753883923853646),
    @(7, "FAPs", "Thbs2", "Itga6", "FAPs", 3, 1, 188.2309416666667, 564.692825, 0.9663705600283187, 0.9663705600283184, 3, 1, 0.4209206666666667, 1.262762, 0.001702901995496819, 0.001702901995496819, 79.2302934536278, 713.0726410826501, 0.001645634355061603, 0.001645634355061602),
    @(8, "FAPs", "Thbs2", "Itga6", "M2", 3, 1, 188.2309416666667, 564.692825, 0.9663705600283187, 0.9663705600283184, 3, 1, 45.70525533333333, 137.115766, 0.184907933193646, 0.184907933193646, 8603.143250508772, 77428.28925457895, 0.1786895829540227, 0.1786895829540226),
    @(9, "FAPs", "Thbs2", "Itga6", "sCs", 3, 1, 188.2309416666667, 564.692825, 0.9663705600283187, 0.9663705600283184, 3, 1, 8.223696, 24.671088, 0.0332702797409786, 0.0332702797409786, 1547.9540420604, 13931.5863785436, 0.03215141886558832, 0.0321514188655883),
    @(10, "M2", "Thbs2", "Itga6", "ECs", 1, 0.3333333333333333, 0.08498099999999999, 0.254943, 0.000436289251037145, 0.000436289251037145, 3, 1, 192.8285726666667, 578.485718, 0.7801188850698786, 0.7801188850698786, 16.386764933786, 147.480884404074, 0.00034035748408707, 0.0003403574840870699),
    @(11, "M2", "Thbs2", "Itga6", "FAPs", 1, 0.3333333333333333, 0.08498099999999999, 0.254943, 0.000436289251037145, 0.000436289251037145, 3, 1, 0.4209206666666667, 1.262762, 0.001702901995496819, 0.001702901995496819, 0.035770259174, 0.321932332566, 0.000000742957836204966900696043, 0.000000742957836204966688937806),
    @(12, "M2", "Thbs2", "Itga6", "M2", 1, 0.3333333333333333, 0.08498099999999999, 0.254943, 0.000436289251037145, 0.000436289251037145, 3, 1, 45.70525533333333, 137.115766, 0.184907933193646, 0.184907933193646, 3.884078303481999, 34.956704731338, 0.000080673343683882275359682545, 0.000080673343683882275359682545),
    @(13, "M2", "Thbs2", "Itga6", "sCs", 1, 0.3333333333333333, 0.08498099999999999, 0.254943, 0.000436289251037145, 0.000436289251037145, 3, 1, 8.223696, 24.671088, 0.0332702797409786, 0.0332702797409786, 0.6988579097759999, 6.289721187984, 0.00001451546542998786068468485, 0.000014515465429987850520289483),
    @(14, "sCs", "Thbs2", "Itga6", "ECs", 3, 1, 4.541390666666667, 13.624172, 0.02331532851610455, 0.02331532851610455, 3, 1, 192.8285726666667, 578.485718, 0.7801188850698786, 0.7801188850698786, 875.7098801750551, 7881.388921575496, 0.01818872808702143, 0.01818872808702143),
    @(15, "sCs", "Thbs2", "Itga6", "FAPs", 3, 1, 4.541390666666667, 13.624172, 0.02331532851610455, 0.02331532851610455, 3, 1, 0.4209206666666667, 1.262762, 0.001702901995496819, 0.001702901995496819, 1.911565187007111, 17.204086683064, 0.000039703719455738327936359516, 0.000039703719455738321160095938),
    @(16, "sCs", "Thbs2", "Itga6", "M2", 3, 1, 4.541390666666667, 13.624172, 0.02331532851610455, 0.02331532851610455, 3, 1, 45.70525533333333, 137.115766, 0.184907933193646, 0.184907933193646, 207.5654199884169, 1868.088779895752, 0.004311189207643771, 0.004311189207643771),
    @(17, "sCs", "Thbs2", "Itga6", "sCs", 3, 1, 4.541390666666667, 13.624172, 0.02331532851610455, 0.02331532851610455, 3, 1, 8.223696, 24.671088, 0.0332702797409786, 0.0332702797409786, 37.347016259904, 336.123146339136, 0.0007757075019836141, 0.0007757075019836138),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    for ($c = 5; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c]
    }
}
